# Commit: "Fixed up Myxicola and a few more loose ends in second review"
#
# On the Materials sheet:
#  - drop the Taxon_Local_ID, suborder, infraorder and superfamily columns
#    (they don't belong on this Darwin-Core-flavoured sheet)
#  - fix the scientificNameAuthorship mapping placeholder
#    (${summary.Author} -> ${summary.authority})

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

$headerRow = 1
$lastCol = $ws.UsedRange.Columns.Count

# Column headers we need to remove entirely from the Materials sheet.
$namesToDelete = @("Taxon_Local_ID", "suborder", "infraorder", "superfamily")

# Collect the column indices for the headers we want to drop.
$colsToDelete = @()
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($headerRow, $c).Value()
    if ($namesToDelete -contains $header) {
        $colsToDelete += $c
    }
}

# Delete from right to left so earlier indices stay valid.
$colsToDelete = $colsToDelete | Sort-Object -Descending
foreach ($c in $colsToDelete) {
    $ws.Cells.Item($headerRow, $c).EntireColumn.Delete()
}

# Fix the authorship placeholder wherever it currently lives.
$newLastCol = $ws.UsedRange.Columns.Count
for ($c = 1; $c -le $newLastCol; $c++) {
    $header = $ws.Cells.Item($headerRow, $c).Value()
    if ($header -eq "scientificNameAuthorship") {
        $ws.Cells.Item(2, $c).Value = "`${summary.authority}"
    }
}
